# "Filled in my execution for the sprint"
# Row 16 (User Story U9 / Task T13) gets "M:1" entered for each of the
# five sprint days (Day 1 .. Day 5 -> columns E..I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16:I16").Value = "M:1"

# Leave the view the way the author left it: scrolled down a few rows with
# the last-filled cell (I16) selected.
$ws.Range("I16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
